$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("sdfsdf", "gsdfgd", " ", "df", "fsd", " ", "444", "333"),
    @("25 f4 Address", "25 f4 city", " ", "25 f4 first", "25 f4 last", " ", "2544", "{{address}}"),
    @("address f5", "city f5", " ", "first f5", "last f5", " ", "555", "5555"),
    @("f3 address", "f2 city", " ", "f2. first", "f2 last", " ", "222", "2222"),
    @("Adress Flow 3 ( Step 2)", "City Flow 3 ( Step 5@)", " ", "First Flow 3", "last Flow 3 ( Step 5)", " ", "333", "3333"),
    @("f3 second address", "f3 4 city", " ", "f3 first (s1)", "f3 five last", " ", "333", "3333"),
    @("hhh", "hhh", " ", "gg", "hhh", " ", "777", "777"),
    @("25 f2 Address (2)", "25 f2 City ( 3)", " ", "25 f2 First (1)", "25 F2 Last (4)", "2502834949444", "252", "250205"),
    @("goo", "ho chi minh", " ", "firs", "last", " ", "5566", "66"),
    @("Adrress 25 f1 @#$%^&*!(#)#*", "25 F1 City", " ", "25 f1 first name", "25 f1 last name", "25012334567955", "251", "2501"),
    @("Address Flow 4", "City Flow 4", " ", "Flow 4 First", "Flow 4 Last", " ", "444", "4444"),
    @("Flow 6 Addresst", "Flow 6 City", " ", "Flow 6 First Name", "Flow 7 Lasst Name", " ", "666", "6666"),
    @("Address Flow 1", "City Flow 1", " ", "first flow 1", "last flow 2", " ", "111", "1111"),
    @("sdafsd", "dgdfg", " ", "first", "ggg", "095756756757656", "444", "fsdfsd"),
    @("f5 address 25", "f5 city 25", " ", "f5 first 25", "f5 last 25", "{{credit_debit_number}}", "2505", "250505"),
    @("Adress flow 2", "City Flow 2", " ", "first", "last flow 2", " ", "222", "2222")
)

$startRow = 18
$endRow = $startRow + $newRows.Count - 1

# Format the whole target range as Text first so purely-numeric looking
# entries (e.g. "444", "095756756757656") are stored as strings, matching
# the rest of the sheet, instead of being auto-coerced into numbers.
$fillRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 8))
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Revert to the plain/default style (matching the rest of the sheet's
# unstyled data rows) now that the text values are safely stored; only the
# Text number format was needed transiently to stop numeric-looking strings
# from being coerced into real numbers.
$fillRange.Style = "Normal"
